$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "value" header (shared string) to "first_release_value"
$ws.Range("B1").Value = "first_release_value"

# Extend the date-formatted style (currently only on A2) down through A22
# before touching values, so every new date cell inherits style index 2
# (the custom date number format + font/border/alignment) from A2.
$ws.Range("A2").Copy($ws.Range("A3:A22"))

# First row's date changes from 39400 to 38717, and its paired value is
# removed entirely (the value series shifts down one row).
$ws.Range("A2").Value = 38717
$ws.Range("B2").ClearContents()

$dates = @(39082, 39447, 39813, 40178, 40543, 40908, 41274, 41639, 42004, 42369, 42735, 43100, 43465, 43830, 44196, 44561, 44926, 45291, 45657, 46022)
$values = @(0.8557439673732903, 1.788430953138542, 1.807765267947059, 1.267704211901699, 0.815159612280536, 2.321967807433256, 1.090188641041823, 0.0688236519329477, 1.078691045907165, 2.724993727165903, 3.366078187926935, 1.625351934832997, 0.8215983724355613, 2.013614902241634, 2.422659863072885, 2.094164808525223, 0.6457774251209525, 1.355088394421644, 0.9622045007620983)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    if ($i -lt $values.Length) {
        $ws.Cells.Item($row, 2).Value = $values[$i]
    }
}
